$d = $word.ActiveDocument

# --- Paragraph 1: "This is a Microsoft word document." -> append trailing
# spaces to the existing run, then append a new red-colored parenthetical
# note, inserted as three separate runs (mirroring the target XML's run
# boundaries).

$p1 = $d.Paragraphs(1).Range

# Collapsed insertion point right before the paragraph mark (End-1, since
# the paragraph Range includes the trailing paragraph-mark character).
$ins = $d.Range($p1.End - 1, $p1.End - 1)

# Two trailing spaces, staying in the original (unformatted) run.
$ins.InsertAfter("  ")
$ins = $d.Range($ins.End, $ins.End)

# First red run.
$start = $ins.Start
$ins.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$d.Range($start, $ins.End).Font.Color = 255
$ins = $d.Range($ins.End, $ins.End)

# Second red run.
$start = $ins.Start
$ins.InsertAfter("rsion for main branch")
$d.Range($start, $ins.End).Font.Color = 255
$ins = $d.Range($ins.End, $ins.End)

# Third red run.
$start = $ins.Start
$ins.InsertAfter(")")
$d.Range($start, $ins.End).Font.Color = 255

# --- Remove the final paragraph ("...ank God almighty, we are free at
# last.") entirely, merging it away (its pPr/rPr NormalWeb styling goes
# away with it).
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastParaIndex)
$lastPara.Range.Delete()

# --- Prune the unused styles that were left over from copy/pasted web
# content (Heading 2/4 + their linked char styles, Hyperlink, and various
# custom "subscribe"/"audio-tool"/"podcast" character & paragraph styles
# pulled in from the source page). None of these are referenced by any
# paragraph in the document body, so removing them is purely a cleanup.
#
# NOTE: styles must be deleted in reverse definition order - deleting an
# earlier-defined style first shifts the collection's indices and this
# host's by-name Styles() lookup resolves stale indices for later names,
# which crashes the interpreter. Walking back-to-front keeps every
# not-yet-deleted style's original index stable.
$stylesToRemove = @(
  "podcast-toolssubscribe-links",
  "generic-title",
  "subscribe-more-info",
  "subscribe",
  "audio-tool",
  "Heading4Char",
  "Heading2Char",
  "Hyperlink",
  "apple-converted-space",
  "Heading4",
  "Heading2"
)

foreach ($styleId in $stylesToRemove) {
  $style = $d.Styles($styleId)
  $style.Delete()
}
